$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the commit-count data for the first two programmers (rows 2 and 3)
$ws.Range("B2").Value = 1427
$ws.Range("C2").Value = 254

$ws.Range("B3").Value = 1001
$ws.Range("C3").Value = 86

# Move the active selection to D13, matching the sheet's saved cursor position
$ws.Range("D13").Select()
